$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 123594.516
$ws.Range("J17").Value = 125907.625
$ws.Range("L17").Value = 377722.875
$ws.Range("N17").Value = -378058.875

$ws.Range("H18").Value = 3313.8333
$ws.Range("I18").Value = 2220.75
$ws.Range("K18").Value = 2220.75
$ws.Range("M18").Value = -1936.75

$ws.Range("H62").Value = 2091
$ws.Range("J62").Value = 2205.2856
$ws.Range("L62").Value = 2205.2856
$ws.Range("N62").Value = -3453.2856

$ws.Range("H65").Value = 2091
$ws.Range("J65").Value = 2205.2856
$ws.Range("L65").Value = 11026.428
$ws.Range("N65").Value = -17266.428

$ws.Range("H80").Value = 12500814
$ws.Range("I80").Value = 577.1667
$ws.Range("J80").Value = 31251168
$ws.Range("K80").Value = 1731.5001
$ws.Range("L80").Value = 93753504
$ws.Range("M80").Value = -733.5001
$ws.Range("N80").Value = -93755500

$ws.Range("H83").Value = 12500814
$ws.Range("I83").Value = 577.1667
$ws.Range("J83").Value = 31251168
$ws.Range("K83").Value = 5194.5003
$ws.Range("L83").Value = 281260512
$ws.Range("M83").Value = -202.5002999999997
$ws.Range("N83").Value = -281270496

$ws.Range("H86").Value = 50003670
$ws.Range("I86").Value = 71432280
$ws.Range("J86").Value = 3582.5
$ws.Range("K86").Value = 71432280
$ws.Range("L86").Value = 3582.5
$ws.Range("M86").Value = -71431157
$ws.Range("N86").Value = -5828.5

$ws.Range("H88").Value = 3925.75
$ws.Range("I88").Value = 1368.75
$ws.Range("J88").Value = 5204.25
$ws.Range("K88").Value = 1368.75
$ws.Range("L88").Value = 5204.25
$ws.Range("M88").Value = -962.75
$ws.Range("N88").Value = -6016.25

$ws.Range("H89").Value = 50003670
$ws.Range("I89").Value = 71432280
$ws.Range("J89").Value = 3582.5
$ws.Range("K89").Value = 357161400
$ws.Range("L89").Value = 17912.5
$ws.Range("M89").Value = -357155784
$ws.Range("N89").Value = -29144.5

$ws.Range("H91").Value = 3925.75
$ws.Range("I91").Value = 1368.75
$ws.Range("J91").Value = 5204.25
$ws.Range("K91").Value = 1368.75
$ws.Range("L91").Value = 5204.25
$ws.Range("M91").Value = 35.25
$ws.Range("N91").Value = -8012.25

$ws.Range("H98").Value = 3351.639
$ws.Range("I98").Value = 1890.75
$ws.Range("K98").Value = 1890.75
$ws.Range("M98").Value = -392.75

$ws.Range("H122").Value = 3351.639
$ws.Range("I122").Value = 1890.75
$ws.Range("K122").Value = 5672.25
$ws.Range("M122").Value = -3222.25

$ws.Range("H131").Value = 7417
$ws.Range("I131").Value = 7772.5
$ws.Range("K131").Value = 23317.5
$ws.Range("M131").Value = -18277.5

$ws.Range("H132").Value = 54141.844
$ws.Range("I132").Value = 60276.176
$ws.Range("K132").Value = 180828.528
$ws.Range("M132").Value = -178298.528

$ws.Range("H137").Value = 2778807.2
$ws.Range("I137").Value = 1030.3928
$ws.Range("J137").Value = 12501026
$ws.Range("K137").Value = 3091.1784
$ws.Range("L137").Value = 37503078
$ws.Range("M137").Value = -541.1784000000002
$ws.Range("N137").Value = -37508178

$ws.Range("H138").Value = 5258.712
$ws.Range("J138").Value = 3337.739
$ws.Range("L138").Value = 10013.217
$ws.Range("N138").Value = -20293.217

$ws.Range("H141").Value = 1303.0667
$ws.Range("I141").Value = 1039
$ws.Range("K141").Value = 3117
$ws.Range("M141").Value = 2063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 44998
$ws.Range("I37").Value = 44998
$ws.Range("K37").Value = 44998
$ws.Range("M37").Value = -44725

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H61").Value = 1523948
$ws.Range("I61").Value = 35848.97
$ws.Range("K61").Value = 35848.97
$ws.Range("M61").Value = -35636.97

$ws.Range("H74").Value = 380151.94
$ws.Range("I74").Value = 1923.375
$ws.Range("J74").Value = 828422.8
$ws.Range("K74").Value = 1923.375
$ws.Range("L74").Value = 828422.8
$ws.Range("M74").Value = -1049.375
$ws.Range("N74").Value = -830170.8

$ws.Range("H77").Value = 380151.94
$ws.Range("I77").Value = 1923.375
$ws.Range("J77").Value = 828422.8
$ws.Range("K77").Value = 9616.875
$ws.Range("L77").Value = 4142114
$ws.Range("M77").Value = -5248.875
$ws.Range("N77").Value = -4150850

$ws.Range("H97").Value = 4921.091
$ws.Range("I97").Value = 5635.1055
$ws.Range("K97").Value = 5635.1055
$ws.Range("M97").Value = -5139.1055

$ws.Range("H122").Value = 1559.8
$ws.Range("I122").Value = 1449.75
$ws.Range("K122").Value = 4349.25
$ws.Range("M122").Value = -1899.25

$ws.Range("H132").Value = 2998.878
$ws.Range("I132").Value = 2679.2666
$ws.Range("J132").Value = 3870.5454
$ws.Range("K132").Value = 8037.7998
$ws.Range("L132").Value = 11611.6362
$ws.Range("M132").Value = -5507.7998
$ws.Range("N132").Value = -16671.6362

$ws.Range("H136").Value = 1523948
$ws.Range("I136").Value = 35848.97
$ws.Range("K136").Value = 107546.91
$ws.Range("M136").Value = -104996.91

$ws.Range("H141").Value = 85902.71000000001
$ws.Range("J141").Value = 85902.71000000001
$ws.Range("L141").Value = 85902.71000000001
$ws.Range("N141").Value = -96262.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2084.0715
$ws.Range("I94").Value = 1414.95
$ws.Range("J94").Value = 3756.875
$ws.Range("K94").Value = 1414.95
$ws.Range("L94").Value = 3756.875
$ws.Range("M94").Value = -963.95
$ws.Range("N94").Value = -4658.875

$ws.Range("H105").Value = 7986.3335
$ws.Range("I105").Value = 8858.23
$ws.Range("K105").Value = 8858.23
$ws.Range("M105").Value = -7111.23

$ws.Range("H134").Value = 34618256
$ws.Range("I134").Value = 2844.4
$ws.Range("K134").Value = 8533.200000000001
$ws.Range("M134").Value = -5998.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2003040
$ws.Range("I6").Value = 2003040
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2003040
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2002927
$ws.Range("N6").ClearContents()

$ws.Range("H31").Value = 2255.966
$ws.Range("I31").Value = 1651.6129
$ws.Range("J31").Value = 2925.0715
$ws.Range("K31").Value = 1651.6129
$ws.Range("L31").Value = 2925.0715
$ws.Range("M31").Value = -1356.6129
$ws.Range("N31").Value = -3515.0715

$ws.Range("H34").Value = 2255.966
$ws.Range("I34").Value = 1651.6129
$ws.Range("J34").Value = 2925.0715
$ws.Range("K34").Value = 1651.6129
$ws.Range("L34").Value = 2925.0715
$ws.Range("M34").Value = -1449.6129
$ws.Range("N34").Value = -3329.0715

$ws.Range("H58").Value = 2306.0967
$ws.Range("I58").Value = 2395.8
$ws.Range("J58").Value = 2263.3809
$ws.Range("K58").Value = 2395.8
$ws.Range("L58").Value = 2263.3809
$ws.Range("M58").Value = -2192.8
$ws.Range("N58").Value = -2669.3809

$ws.Range("H62").Value = 6640.9
$ws.Range("I62").Value = 4735.8335
$ws.Range("J62").Value = 9498.5
$ws.Range("K62").Value = 4735.8335
$ws.Range("L62").Value = 9498.5
$ws.Range("M62").Value = -4111.8335
$ws.Range("N62").Value = -10746.5

$ws.Range("H65").Value = 6640.9
$ws.Range("I65").Value = 4735.8335
$ws.Range("J65").Value = 9498.5
$ws.Range("K65").Value = 23679.1675
$ws.Range("L65").Value = 47492.5
$ws.Range("M65").Value = -20559.1675
$ws.Range("N65").Value = -53732.5

$ws.Range("H99").Value = 2224497
$ws.Range("I99").Value = 2859283
$ws.Range("J99").Value = 2746.5
$ws.Range("K99").Value = 2859283
$ws.Range("L99").Value = 2746.5
$ws.Range("M99").Value = -2857785
$ws.Range("N99").Value = -5742.5

$ws.Range("H107").Value = 1720.4482
$ws.Range("I107").Value = 1494.05
$ws.Range("J107").Value = 2223.5557
$ws.Range("K107").Value = 1494.05
$ws.Range("L107").Value = 2223.5557
$ws.Range("M107").Value = 425.95
$ws.Range("N107").Value = -6063.5557

$ws.Range("H126").Value = 2224497
$ws.Range("I126").Value = 2859283
$ws.Range("J126").Value = 2746.5
$ws.Range("K126").Value = 8577849
$ws.Range("L126").Value = 8239.5
$ws.Range("M126").Value = -8575379
$ws.Range("N126").Value = -13179.5

$ws.Range("H132").Value = 24582.25
$ws.Range("I132").Value = 32987.812
$ws.Range("J132").Value = 2167.4167
$ws.Range("K132").Value = 98963.43599999999
$ws.Range("L132").Value = 6502.250100000001
$ws.Range("M132").Value = -96433.43599999999
$ws.Range("N132").Value = -11562.2501

$ws.Range("H134").Value = 2463.238
$ws.Range("I134").Value = 1955.0769
$ws.Range("J134").Value = 3289
$ws.Range("K134").Value = 5865.2307
$ws.Range("L134").Value = 9867
$ws.Range("M134").Value = -3330.2307
$ws.Range("N134").Value = -14937

$ws.Range("H136").Value = 2306.0967
$ws.Range("I136").Value = 2395.8
$ws.Range("J136").Value = 2263.3809
$ws.Range("K136").Value = 7187.400000000001
$ws.Range("L136").Value = 6790.1427
$ws.Range("M136").Value = -4637.400000000001
$ws.Range("N136").Value = -11890.1427

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 669.8570999999999
$ws.Range("I5").Value = 466.46667
$ws.Range("J5").Value = 1178.3334
$ws.Range("K5").Value = 1399.40001
$ws.Range("L5").Value = 3535.0002
$ws.Range("M5").Value = -1287.40001
$ws.Range("N5").Value = -3759.0002

$ws.Range("H15").Value = 2230
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 2966.6667
$ws.Range("K15").Value = 60
$ws.Range("L15").Value = 8900.000100000001
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = -9180.000100000001

$ws.Range("H17").Value = 953.3125
$ws.Range("I17").Value = 263.5
$ws.Range("J17").Value = 1183.25
$ws.Range("K17").Value = 790.5
$ws.Range("L17").Value = 3549.75
$ws.Range("M17").Value = -621.5
$ws.Range("N17").Value = -3887.75

$ws.Range("H21").Value = 281.27777
$ws.Range("I21").Value = 209.58824
$ws.Range("J21").Value = 1500
$ws.Range("K21").Value = 628.76472
$ws.Range("L21").Value = 4500
$ws.Range("M21").Value = -455.76472
$ws.Range("N21").Value = -4846

$ws.Range("H26").Value = 344.84616
$ws.Range("I26").Value = 475
$ws.Range("J26").Value = 233.28572
$ws.Range("K26").Value = 1425
$ws.Range("L26").Value = 699.85716
$ws.Range("M26").Value = -1137
$ws.Range("N26").Value = -1275.85716

$ws.Range("H32").Value = 2966.3333
$ws.Range("I32").Value = 2966
$ws.Range("J32").Value = 2966.6667
$ws.Range("K32").Value = 8898
$ws.Range("L32").Value = 8900.000100000001
$ws.Range("M32").Value = -8615
$ws.Range("N32").Value = -9466.000100000001

$ws.Range("H36").Value = 278
$ws.Range("I36").Value = 487
$ws.Range("J36").Value = 69
$ws.Range("K36").Value = 1461
$ws.Range("L36").Value = 207
$ws.Range("M36").Value = -1292
$ws.Range("N36").Value = -545

$ws.Range("H39").Value = 6560.7144
$ws.Range("I39").Value = 855
$ws.Range("K39").Value = 2565
$ws.Range("M39").Value = -2271

$ws.Range("H41").Value = 1112.8572
$ws.Range("I41").Value = 1050
$ws.Range("K41").Value = 3150
$ws.Range("M41").Value = -2812

$ws.Range("H42").Value = 33350900
$ws.Range("J42").Value = 33350900
$ws.Range("L42").Value = 100052700
$ws.Range("N42").Value = -100053768

$ws.Range("H49").Value = 3501.5
$ws.Range("I49").Value = 3501.5
$ws.Range("K49").Value = 10504.5
$ws.Range("M49").Value = -10348.5

$ws.Range("H50").Value = 1791.5
$ws.Range("I50").Value = 253.8
$ws.Range("J50").Value = 3329.2
$ws.Range("K50").Value = 761.4000000000001
$ws.Range("L50").Value = 9987.599999999999
$ws.Range("M50").Value = -280.4000000000001
$ws.Range("N50").Value = -10949.6

$ws.Range("H53").Value = 1791.5
$ws.Range("I53").Value = 253.8
$ws.Range("J53").Value = 3329.2
$ws.Range("K53").Value = 761.4000000000001
$ws.Range("L53").Value = 9987.599999999999
$ws.Range("M53").Value = -280.4000000000001
$ws.Range("N53").Value = -10949.6

$ws.Range("H69").Value = 7803.75
$ws.Range("J69").Value = 5997
$ws.Range("L69").Value = 17991
$ws.Range("N69").Value = -19613

$ws.Range("H72").Value = 7803.75
$ws.Range("J72").Value = 5997
$ws.Range("L72").Value = 53973
$ws.Range("N72").Value = -62085

$ws.Range("H102").Value = 17874
$ws.Range("J102").Value = 17874
$ws.Range("L102").Value = 53622
$ws.Range("N102").Value = -58490

$ws.Range("H114").Value = 3893.1177
$ws.Range("I114").Value = 1180.125
$ws.Range("J114").Value = 6304.6665
$ws.Range("K114").Value = 3540.375
$ws.Range("L114").Value = 18913.9995
$ws.Range("M114").Value = -286.375
$ws.Range("N114").Value = -25421.9995

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H129").Value = 8911.579
$ws.Range("I129").Value = 983.6923
$ws.Range("J129").Value = 26088.666
$ws.Range("K129").Value = 2951.0769
$ws.Range("L129").Value = 78265.99800000001
$ws.Range("M129").Value = 2048.9231
$ws.Range("N129").Value = -88265.99800000001

$ws.Range("H131").Value = 8268090.5
$ws.Range("I131").Value = 22730410
$ws.Range("J131").Value = 3907.8572
$ws.Range("K131").Value = 68191230
$ws.Range("L131").Value = 11723.5716
$ws.Range("M131").Value = -68186190
$ws.Range("N131").Value = -21803.5716

$ws.Range("H134").Value = 4553.3076
$ws.Range("I134").Value = 2653.9092
$ws.Range("K134").Value = 7961.7276
$ws.Range("M134").Value = -2891.7276

$ws.Range("H135").Value = 669.8570999999999
$ws.Range("I135").Value = 466.46667
$ws.Range("J135").Value = 1178.3334
$ws.Range("K135").Value = 4198.20003
$ws.Range("L135").Value = 10605.0006
$ws.Range("M135").Value = -1663.20003
$ws.Range("N135").Value = -15675.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H122").Value = 1656.0444
$ws.Range("I122").Value = 1625.5312
$ws.Range("K122").Value = 4876.5936
$ws.Range("M122").Value = -2426.5936

$ws.Range("H126").Value = 2099
$ws.Range("I126").Value = 2165
$ws.Range("K126").Value = 6495
$ws.Range("M126").Value = -4025

$ws.Range("H132").Value = 734404.0600000001
$ws.Range("I132").Value = 1583.5238
$ws.Range("J132").Value = 1833634.9
$ws.Range("K132").Value = 4750.5714
$ws.Range("L132").Value = 5500904.699999999
$ws.Range("M132").Value = -2220.5714
$ws.Range("N132").Value = -5505964.699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6630.615
$ws.Range("I7").Value = 2407.1428
$ws.Range("K7").Value = 2407.1428
$ws.Range("M7").Value = -2295.1428

$ws.Range("H14").Value = 82423.75
$ws.Range("I14").Value = 106566.664
$ws.Range("K14").Value = 106566.664
$ws.Range("M14").Value = -106394.664

$ws.Range("H23").Value = 1255846.5
$ws.Range("I23").Value = 1255846.5
$ws.Range("K23").Value = 1255846.5
$ws.Range("M23").Value = -1255616.5

$ws.Range("H36").Value = 65000
$ws.Range("J36").Value = 65000
$ws.Range("L36").Value = 65000
$ws.Range("N36").Value = -66124

$ws.Range("H40").Value = 1963.2858
$ws.Range("I40").Value = 1142
$ws.Range("J40").Value = 2784.5715
$ws.Range("K40").Value = 1142
$ws.Range("L40").Value = 2784.5715
$ws.Range("M40").Value = -1006
$ws.Range("N40").Value = -3056.5715

$ws.Range("H68").Value = 4193.222
$ws.Range("I68").Value = 5209.4
$ws.Range("K68").Value = 5209.4
$ws.Range("M68").Value = -4460.4

$ws.Range("H71").Value = 4193.222
$ws.Range("I71").Value = 5209.4
$ws.Range("K71").Value = 26047
$ws.Range("M71").Value = -22303

$ws.Range("H93").Value = 1227.6666
$ws.Range("I93").Value = 1227.6666
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1227.6666
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 20.33339999999998
$ws.Range("N93").ClearContents()

$ws.Range("H101").Value = 29816
$ws.Range("J101").Value = 29816
$ws.Range("L101").Value = 29816
$ws.Range("N101").Value = -36306

$ws.Range("H122").Value = 2300.524
$ws.Range("I122").Value = 2088.76
$ws.Range("K122").Value = 6266.280000000001
$ws.Range("M122").Value = -3816.280000000001

$ws.Range("H125").Value = 113972.5
$ws.Range("J125").Value = 113972.5
$ws.Range("L125").Value = 113972.5
$ws.Range("N125").Value = -123812.5

$ws.Range("H126").Value = 6630.615
$ws.Range("I126").Value = 2407.1428
$ws.Range("K126").Value = 7221.428400000001
$ws.Range("M126").Value = -4751.428400000001

$ws.Range("H132").Value = 2618.7917
$ws.Range("I132").Value = 2208.6333
$ws.Range("K132").Value = 6625.8999
$ws.Range("M132").Value = -4095.8999

$ws.Range("H136").Value = 2064.7036
$ws.Range("I136").Value = 2716.7646
$ws.Range("K136").Value = 8150.293799999999
$ws.Range("M136").Value = -5600.293799999999

$ws.Range("H140").Value = 84994.5
$ws.Range("J140").Value = 84994.5
$ws.Range("L140").Value = 84994.5
$ws.Range("N140").Value = -95354.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 10000
$ws.Range("K24").Value = 10000
$ws.Range("M24").Value = -9770

$ws.Range("H81").Value = 65399.688
$ws.Range("I81").Value = 3087.5
$ws.Range("K81").Value = 6175
$ws.Range("M81").Value = -5114

$ws.Range("H84").Value = 65399.688
$ws.Range("I84").Value = 3087.5
$ws.Range("K84").Value = 30875
$ws.Range("M84").Value = -25571

$ws.Range("H100").Value = 888.375
$ws.Range("I100").Value = 886.7143
$ws.Range("K100").Value = 1773.4286
$ws.Range("M100").Value = -1232.4286

$ws.Range("H107").Value = 497.14285
$ws.Range("I107").Value = 504.125
$ws.Range("K107").Value = 1512.375
$ws.Range("M107").Value = 407.625

$ws.Range("H122").Value = 2221.9
$ws.Range("I122").Value = 1874.0588
$ws.Range("K122").Value = 5622.1764
$ws.Range("M122").Value = -3172.1764

$ws.Range("H126").Value = 1929.2106
$ws.Range("I126").Value = 1638.5
$ws.Range("K126").Value = 4915.5
$ws.Range("M126").Value = -2445.5

$ws.Range("H132").Value = 1883.5
$ws.Range("I132").Value = 1568.1555
$ws.Range("J132").Value = 2409.074
$ws.Range("K132").Value = 4704.4665
$ws.Range("L132").Value = 7227.222
$ws.Range("M132").Value = -2174.4665
$ws.Range("N132").Value = -12287.222

$ws.Range("H136").Value = 33313.742
$ws.Range("J136").Value = 2719
$ws.Range("L136").Value = 8157
$ws.Range("N136").Value = -13257
